# Mise à jour de l'application
# Adds a new date column (CW) to the "Présences" attendance sheet, one day
# after the existing last column (CV = 06/01/2026), i.e. CW = 07/01/2026
# (Excel serial 46029), and fills in the attendance marks for that date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell CW1: new date, style copied from CV1 (date format, centered) ---
$ws.Range("CV1").Copy()
$ws.Range("CW1").PasteSpecial(-4122)
$ws.Range("CW1").Value = 46029

# --- Data cells CW2:CW30: attendance marks for the new date ---
# (rows 12, 21 and 23 have no entry for this date, matching the existing
# gaps already present for those players/rows)
$ws.Range("CW2").Value = "P"
$ws.Range("CW2").HorizontalAlignment = -4108

$ws.Range("CW3").Value = "P"
$ws.Range("CW3").HorizontalAlignment = -4108

$ws.Range("CW4").Value = "P"
$ws.Range("CW4").HorizontalAlignment = -4108

$ws.Range("CW5").Value = "P"
$ws.Range("CW5").HorizontalAlignment = -4108

$ws.Range("CW6").Value = "P"
$ws.Range("CW6").HorizontalAlignment = -4108

$ws.Range("CW7").Value = "P"
$ws.Range("CW7").HorizontalAlignment = -4108

$ws.Range("CW8").Value = "P"
$ws.Range("CW8").HorizontalAlignment = -4108

$ws.Range("CW9").Value = "RH"
$ws.Range("CW9").HorizontalAlignment = -4108

$ws.Range("CW10").Value = "P"
$ws.Range("CW10").HorizontalAlignment = -4108

$ws.Range("CW11").Value = "P"
$ws.Range("CW11").HorizontalAlignment = -4108

$ws.Range("CW13").Value = "B"
$ws.Range("CW13").HorizontalAlignment = -4108

$ws.Range("CW14").Value = "P"
$ws.Range("CW14").HorizontalAlignment = -4108

$ws.Range("CW15").Value = "P"
$ws.Range("CW15").HorizontalAlignment = -4108

$ws.Range("CW16").Value = "REP"
$ws.Range("CW16").HorizontalAlignment = -4108

$ws.Range("CW17").Value = "P"
$ws.Range("CW17").HorizontalAlignment = -4108

$ws.Range("CW18").Value = "P"
$ws.Range("CW18").HorizontalAlignment = -4108

$ws.Range("CW19").Value = "P"
$ws.Range("CW19").HorizontalAlignment = -4108

$ws.Range("CW20").Value = "P"
$ws.Range("CW20").HorizontalAlignment = -4108

$ws.Range("CW22").Value = "P"
$ws.Range("CW22").HorizontalAlignment = -4108

$ws.Range("CW24").Value = "P"
$ws.Range("CW24").HorizontalAlignment = -4108

$ws.Range("CW25").Value = "A"
$ws.Range("CW25").HorizontalAlignment = -4108

$ws.Range("CW26").Value = "P"
$ws.Range("CW26").HorizontalAlignment = -4108

$ws.Range("CW27").Value = "P"
$ws.Range("CW27").HorizontalAlignment = -4108

$ws.Range("CW28").Value = "P"
$ws.Range("CW28").HorizontalAlignment = -4108

$ws.Range("CW29").Value = "P"
$ws.Range("CW29").HorizontalAlignment = -4108

$ws.Range("CW30").Value = "P"
$ws.Range("CW30").HorizontalAlignment = -4108

# --- View state: move the selection to the new last column, row 27 ---
$ws.Range("CW27").Select()
